$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) from row 2 to row 271: 45202 -> 45203
$ws.Range("C2:C271").Value = 45203
